# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update worker identification (row 16: Doc #, Name, Period)
$ws.Range("C16").Value = "1050960222"
$ws.Range("D16").Value = "MARIA ALEJANDRA TORRES ESPINOSA"
$ws.Range("E16").Value = "2507"

# Update mora period value and figures
$ws.Range("E11").Value = 1898
$ws.Range("F16").Value = 1898
$ws.Range("G16").Value = 1423500

# Let column D auto-fit to the new (longer) worker name
# (closest reproducible value to the target stored width of 35.453125)
$ws.Columns.Item(4).ColumnWidth = 34.666666666666664
